$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row relabeling ---
$ws.Range("B1").Value = "Diversity metric"
$ws.Range("D1").Value = "% Difference"
$ws.Range("E1").Value = "% Error"

# --- Data: Slope (C), % Difference (D), % Error (E) per model/metric ---
# Lognormal
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1.5
$ws.Range("E2").Value = 1.5

$ws.Range("C3").Value = -0.48
$ws.Range("D3").Value = 42
$ws.Range("E3").Value = 53

$ws.Range("C4").Value = 0.1
$ws.Range("D4").Value = 23
$ws.Range("E4").Value = 20

# Zipf
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0.28000000000000003
$ws.Range("E5").Value = 0.28000000000000003

$ws.Range("C6").Value = -0.53
$ws.Range("D6").Value = 53
$ws.Range("E6").Value = 72

$ws.Range("C7").Value = 0.086
$ws.Range("D7").Value = 41
$ws.Range("E7").Value = 34

# Log-series
$ws.Range("C8").Value = 0.86
$ws.Range("D8").Value = 16
$ws.Range("E8").Value = 14

$ws.Range("C9").Value = -0.16
$ws.Range("D9").Value = 66
$ws.Range("E9").Value = 50

$ws.Range("C10").Value = 0.048
$ws.Range("D10").Value = 92
$ws.Range("E10").Value = 63

# Broken-stick
$ws.Range("C11").Value = 0.73
$ws.Range("D11").Value = 32
$ws.Range("E11").Value = 27

$ws.Range("C12").Value = -0.022
$ws.Range("D12").Value = 170
$ws.Range("E12").Value = 93

$ws.Range("C13").Value = 0.014
$ws.Range("D13").Value = 160
$ws.Range("E13").Value = 89

# --- Formatting: row 1 reverts to auto height (no explicit ht) ---
$ws.Rows.Item(1).AutoFit()

# --- Selection moves from D2 to E2 ---
[void]$ws.Range("E2").Select()
